$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 392.25
$ws.Range("I9").Value = 325
$ws.Range("K9").Value = 325
$ws.Range("M9").Value = -156
$ws.Range("H17").Value = 360.6216
$ws.Range("J17").Value = 306.94287
$ws.Range("L17").Value = 920.82861
$ws.Range("N17").Value = -1256.82861
$ws.Range("H41").Value = 20346
$ws.Range("I41").Value = 432
$ws.Range("K41").Value = 432
$ws.Range("M41").Value = 8
$ws.Range("H62").Value = 4097.3335
$ws.Range("I62").Value = 4371
$ws.Range("J62").Value = 3550
$ws.Range("K62").Value = 4371
$ws.Range("L62").Value = 3550
$ws.Range("M62").Value = -3747
$ws.Range("N62").Value = -4798
$ws.Range("H65").Value = 4097.3335
$ws.Range("I65").Value = 4371
$ws.Range("J65").Value = 3550
$ws.Range("K65").Value = 21855
$ws.Range("L65").Value = 17750
$ws.Range("M65").Value = -18735
$ws.Range("N65").Value = -23990
$ws.Range("H69").Value = 23015.867
$ws.Range("I69").Value = 6066.6665
$ws.Range("J69").Value = 34315.332
$ws.Range("K69").Value = 18199.9995
$ws.Range("L69").Value = 102945.996
$ws.Range("M69").Value = -17325.9995
$ws.Range("N69").Value = -104693.996
$ws.Range("H70").Value = 1316.6666
$ws.Range("I70").Value = 1050
$ws.Range("K70").Value = 3150
$ws.Range("M70").Value = -2880
$ws.Range("H72").Value = 23015.867
$ws.Range("I72").Value = 6066.6665
$ws.Range("J72").Value = 34315.332
$ws.Range("K72").Value = 54599.9985
$ws.Range("L72").Value = 308837.988
$ws.Range("M72").Value = -50231.9985
$ws.Range("N72").Value = -317573.988
$ws.Range("H73").Value = 1316.6666
$ws.Range("I73").Value = 1050
$ws.Range("K73").Value = 3150
$ws.Range("M73").Value = -2214
$ws.Range("H80").Value = 2459.258
$ws.Range("I80").Value = 817.2222
$ws.Range("J80").Value = 3131
$ws.Range("K80").Value = 2451.6666
$ws.Range("L80").Value = 9393
$ws.Range("M80").Value = -1453.6666
$ws.Range("N80").Value = -11389
$ws.Range("H83").Value = 2459.258
$ws.Range("I83").Value = 817.2222
$ws.Range("J83").Value = 3131
$ws.Range("K83").Value = 7354.999800000001
$ws.Range("L83").Value = 28179
$ws.Range("M83").Value = -2362.999800000001
$ws.Range("N83").Value = -38163
$ws.Range("H86").Value = 2504.5
$ws.Range("J86").Value = 2800.5715
$ws.Range("L86").Value = 2800.5715
$ws.Range("N86").Value = -5046.5715
$ws.Range("H89").Value = 2504.5
$ws.Range("J89").Value = 2800.5715
$ws.Range("L89").Value = 14002.8575
$ws.Range("N89").Value = -25234.8575
$ws.Range("H132").Value = 12170.149
$ws.Range("I132").Value = 1186.9762
$ws.Range("K132").Value = 3560.9286
$ws.Range("M132").Value = -1030.9286
$ws.Range("H137").Value = 3554.8147
$ws.Range("I137").Value = 3717.0557
$ws.Range("J137").Value = 3230.3333
$ws.Range("K137").Value = 11151.1671
$ws.Range("L137").Value = 9690.999899999999
$ws.Range("M137").Value = -8601.167099999999
$ws.Range("N137").Value = -14790.9999
$ws.Range("H138").Value = 3104.762
$ws.Range("J138").Value = 4163
$ws.Range("L138").Value = 12489
$ws.Range("N138").Value = -22769

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7802.2354
$ws.Range("I32").Value = 11841.429
$ws.Range("J32").Value = 1277.3846
$ws.Range("K32").Value = 11841.429
$ws.Range("L32").Value = 1277.3846
$ws.Range("M32").Value = -11554.429
$ws.Range("N32").Value = -1851.3846
$ws.Range("H45").Value = 2493.56
$ws.Range("I45").Value = 2080.111
$ws.Range("K45").Value = 2080.111
$ws.Range("M45").Value = -1703.111
$ws.Range("H98").Value = 53083.332
$ws.Range("J98").Value = 53083.332
$ws.Range("L98").Value = 53083.332
$ws.Range("N98").Value = -59073.332
$ws.Range("H103").Value = 37500
$ws.Range("J103").Value = 37500
$ws.Range("L103").Value = 37500
$ws.Range("N103").Value = -39844
$ws.Range("H132").Value = 2353.0667
$ws.Range("I132").Value = 1946.2307
$ws.Range("K132").Value = 5838.6921
$ws.Range("M132").Value = -3308.6921

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H20").Value = 7283.3335
$ws.Range("I20").Value = 4566.6665
$ws.Range("K20").Value = 4566.6665
$ws.Range("M20").Value = -4319.6665
$ws.Range("H86").Value = 1600
$ws.Range("I86").Value = 1200
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -77
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1600
$ws.Range("I89").Value = 1200
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 6000
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -384
$ws.Range("N89").Value = -21232
$ws.Range("H99").Value = 60182.57
$ws.Range("I99").Value = 135270
$ws.Range("K99").Value = 135270
$ws.Range("M99").Value = -133772
$ws.Range("H107").Value = 11578.77
$ws.Range("I107").Value = 12657.625
$ws.Range("K107").Value = 12657.625
$ws.Range("M107").Value = -10737.625
$ws.Range("H134").Value = 864.35
$ws.Range("I134").Value = 864.35
$ws.Range("K134").Value = 2593.05
$ws.Range("M134").Value = -58.05000000000018

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 4596.407
$ws.Range("I107").Value = 737.7059
$ws.Range("J107").Value = 11156.2
$ws.Range("K107").Value = 737.7059
$ws.Range("L107").Value = 11156.2
$ws.Range("M107").Value = 1182.2941
$ws.Range("N107").Value = -14996.2
$ws.Range("H132").Value = 7599
$ws.Range("I132").Value = 5998
$ws.Range("K132").Value = 17994
$ws.Range("M132").Value = -15464
$ws.Range("H141").Value = 55000
$ws.Range("J141").Value = 55000
$ws.Range("L141").Value = 55000
$ws.Range("N141").Value = -65360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8790.6
$ws.Range("I3").Value = 7545.1113
$ws.Range("K3").Value = 22635.3339
$ws.Range("M3").Value = -22523.3339
$ws.Range("H12").Value = 654.26666
$ws.Range("I12").Value = 546.5
$ws.Range("J12").Value = 693.4545000000001
$ws.Range("K12").Value = 1639.5
$ws.Range("L12").Value = 2080.3635
$ws.Range("M12").Value = -1466.5
$ws.Range("N12").Value = -2426.3635
$ws.Range("H15").Value = 163.33333
$ws.Range("I15").Value = 100
$ws.Range("J15").Value = 195
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 585
$ws.Range("M15").Value = -160
$ws.Range("N15").Value = -865
$ws.Range("H39").Value = 1074.7142
$ws.Range("J39").Value = 1124.5
$ws.Range("L39").Value = 3373.5
$ws.Range("N39").Value = -3961.5
$ws.Range("H107").Value = 2479
$ws.Range("I107").Value = 2800.25
$ws.Range("J107").Value = 2362.182
$ws.Range("K107").Value = 8400.75
$ws.Range("L107").Value = 7086.545999999999
$ws.Range("M107").Value = -6480.75
$ws.Range("N107").Value = -10926.546
$ws.Range("H136").Value = 909
$ws.Range("I136").Value = 909
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2727
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 2373
$ws.Range("N136").ClearContents()
$ws.Range("H140").Value = 2895.4546
$ws.Range("I140").Value = 2651.25
$ws.Range("K140").Value = 7953.75
$ws.Range("M140").Value = -2773.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 66361.336
$ws.Range("I80").Value = 187918.5
$ws.Range("J80").Value = 5582.75
$ws.Range("K80").Value = 187918.5
$ws.Range("L80").Value = 5582.75
$ws.Range("M80").Value = -186920.5
$ws.Range("N80").Value = -7578.75
$ws.Range("H83").Value = 66361.336
$ws.Range("I83").Value = 187918.5
$ws.Range("J83").Value = 5582.75
$ws.Range("K83").Value = 939592.5
$ws.Range("L83").Value = 27913.75
$ws.Range("M83").Value = -934600.5
$ws.Range("N83").Value = -37897.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2244.125
$ws.Range("I46").Value = 1293.6
$ws.Range("K46").Value = 1293.6
$ws.Range("M46").Value = -1105.6
$ws.Range("H95").Value = 58500
$ws.Range("J95").Value = 58500
$ws.Range("L95").Value = 58500
$ws.Range("N95").Value = -63992
$ws.Range("H98").Value = 49992.5
$ws.Range("J98").Value = 49992.5
$ws.Range("L98").Value = 49992.5
$ws.Range("N98").Value = -55982.5
$ws.Range("H103").Value = 55900.105
$ws.Range("I103").Value = 35000
$ws.Range("J103").Value = 57061.223
$ws.Range("K103").Value = 35000
$ws.Range("L103").Value = 57061.223
$ws.Range("M103").Value = -33828
$ws.Range("N103").Value = -59405.223
$ws.Range("H136").Value = 2891.423
$ws.Range("I136").Value = 2201.7778
$ws.Range("K136").Value = 6605.3334
$ws.Range("M136").Value = -4055.3334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11114011
$ws.Range("I81").Value = 1999.5
$ws.Range("J81").Value = 18522018
$ws.Range("K81").Value = 3999
$ws.Range("L81").Value = 37044036
$ws.Range("M81").Value = -2938
$ws.Range("N81").Value = -37046158
$ws.Range("H84").Value = 11114011
$ws.Range("I84").Value = 1999.5
$ws.Range("J84").Value = 18522018
$ws.Range("K84").Value = 19995
$ws.Range("L84").Value = 185220180
$ws.Range("M84").Value = -14691
$ws.Range("N84").Value = -185230788
$ws.Range("H136").Value = 1874.5
$ws.Range("I136").Value = 1874.5
$ws.Range("K136").Value = 5623.5
$ws.Range("M136").Value = -3073.5
